$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 414; this shifts the old rows 414-424 down to 417-427
# (formatting, e.g. the date style on column D, is carried down from the row above,
# matching native Excel "Insert" behaviour).
$ws.Rows("414:416").Insert()

# Row 414 - new weekly entry
$ws.Range("A414").Value = 4
$ws.Range("B414").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C414").Value = "Los Lagos"
$ws.Range("D414").Value = 44448
$ws.Range("E414").Value = 10
$ws.Range("F414").Value = 100112002
$ws.Range("G414").Value = "Pimiento"
$ws.Range("H414").Value = "Zafiro rojo"
$ws.Range("I414").Value = "Primera"
$ws.Range("J414").Value = 50
$ws.Range("K414").Value = 44000
$ws.Range("L414").Value = 44000
$ws.Range("M414").Value = 44000
$ws.Range("N414").Value = "$/caja 15 kilos"
$ws.Range("O414").Value = "Región de Arica y Parinacota"
$ws.Range("P414").Value = 2933
$ws.Range("Q414").Value = 15
$ws.Range("R414").Value = "Hortaliza"

# Row 415 - new weekly entry
$ws.Range("A415").Value = 4
$ws.Range("B415").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C415").Value = "Los Lagos"
$ws.Range("D415").Value = 44448
$ws.Range("E415").Value = 10
$ws.Range("F415").Value = 100112002
$ws.Range("G415").Value = "Pimiento"
$ws.Range("H415").Value = "Zafiro verde"
$ws.Range("I415").Value = "Extra"
$ws.Range("J415").Value = 50
$ws.Range("K415").Value = 41000
$ws.Range("L415").Value = 41000
$ws.Range("M415").Value = 41000
$ws.Range("N415").Value = "$/caja 15 kilos"
$ws.Range("O415").Value = "Región de Arica y Parinacota"
$ws.Range("P415").Value = 2733
$ws.Range("Q415").Value = 15
$ws.Range("R415").Value = "Hortaliza"

# Row 416 - new weekly entry
$ws.Range("A416").Value = 4
$ws.Range("B416").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C416").Value = "Los Lagos"
$ws.Range("D416").Value = 44448
$ws.Range("E416").Value = 10
$ws.Range("F416").Value = 100112002
$ws.Range("G416").Value = "Pimiento"
$ws.Range("H416").Value = "Zafiro verde"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 50
$ws.Range("K416").Value = 38000
$ws.Range("L416").Value = 38000
$ws.Range("M416").Value = 38000
$ws.Range("N416").Value = "$/caja 15 kilos"
$ws.Range("O416").Value = "Región de Arica y Parinacota"
$ws.Range("P416").Value = 2533
$ws.Range("Q416").Value = 15
$ws.Range("R416").Value = "Hortaliza"
